$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: clear contents and force Text number format on all cells that will
# be rewritten, so the numeric-looking strings are stored as text (matching
# the original inline-string cell type) instead of being auto-parsed as numbers.
$r = $ws.Range("D2:D18")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("D21:D23")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("D25:D26")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("D40:D50")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("E2:E21")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("E23:E26")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("E39:E50")
$r.ClearContents()
$r.NumberFormat = "@"
$r = $ws.Range("G2:G51")
$r.ClearContents()
$r.NumberFormat = "@"

# Step 2: write the new text values cell by cell.
$ws.Range("D2").Value = "293.72"
$ws.Range("E2").Value = "-2.93%"
$ws.Range("G2").Value = "19"
$ws.Range("D3").Value = "31.06"
$ws.Range("E3").Value = "-2.78%"
$ws.Range("G3").Value = "19"
$ws.Range("D4").Value = "4.886"
$ws.Range("E4").Value = "-2.95%"
$ws.Range("G4").Value = "19"
$ws.Range("D5").Value = "0.07332"
$ws.Range("E5").Value = "-7.13%"
$ws.Range("G5").Value = "19"
$ws.Range("D6").Value = "1.828"
$ws.Range("E6").Value = "-12.38%"
$ws.Range("G6").Value = "19"
$ws.Range("D7").Value = "7.676"
$ws.Range("E7").Value = "-1.91%"
$ws.Range("G7").Value = "19"
$ws.Range("D8").Value = "3.767"
$ws.Range("E8").Value = "-0.83%"
$ws.Range("G8").Value = "19"
$ws.Range("D9").Value = "0.9068"
$ws.Range("E9").Value = "-2.18%"
$ws.Range("G9").Value = "19"
$ws.Range("D10").Value = "0.1653"
$ws.Range("E10").Value = "-6.05%"
$ws.Range("G10").Value = "19"
$ws.Range("D11").Value = "0.07574"
$ws.Range("E11").Value = "-5.41%"
$ws.Range("G11").Value = "19"
$ws.Range("D12").Value = "0.08173"
$ws.Range("E12").Value = "-6.89%"
$ws.Range("G12").Value = "19"
$ws.Range("D13").Value = "0.02992"
$ws.Range("E13").Value = "-3.53%"
$ws.Range("G13").Value = "19"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "-0.28%"
$ws.Range("G14").Value = "19"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").Value = "-0.48%"
$ws.Range("G15").Value = "19"
$ws.Range("D16").Value = "0.005653"
$ws.Range("E16").Value = "-4.68%"
$ws.Range("G16").Value = "19"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("G17").Value = "19"
$ws.Range("D18").Value = "2.097"
$ws.Range("E18").Value = "-8.04%"
$ws.Range("G18").Value = "19"
$ws.Range("E19").Value = "-0.65%"
$ws.Range("G19").Value = "19"
$ws.Range("E20").Value = "1.33%"
$ws.Range("G20").Value = "19"
$ws.Range("D21").Value = "4.369"
$ws.Range("E21").Value = "5.45%"
$ws.Range("G21").Value = "19"
$ws.Range("D22").Value = "0.2002"
$ws.Range("G22").Value = "19"
$ws.Range("D23").Value = "0.04482"
$ws.Range("E23").Value = "-2.59%"
$ws.Range("G23").Value = "19"
$ws.Range("E24").Value = "-0.80%"
$ws.Range("G24").Value = "19"
$ws.Range("D25").Value = "0.004037"
$ws.Range("E25").Value = "-10.48%"
$ws.Range("G25").Value = "19"
$ws.Range("D26").Value = "0.0001252"
$ws.Range("E26").Value = "0.33%"
$ws.Range("G26").Value = "19"
$ws.Range("G27").Value = "19"
$ws.Range("G28").Value = "19"
$ws.Range("G29").Value = "19"
$ws.Range("G30").Value = "19"
$ws.Range("G31").Value = "19"
$ws.Range("G32").Value = "19"
$ws.Range("G33").Value = "19"
$ws.Range("G34").Value = "19"
$ws.Range("G35").Value = "19"
$ws.Range("G36").Value = "19"
$ws.Range("G37").Value = "19"
$ws.Range("G38").Value = "19"
$ws.Range("E39").Value = "-4.85%"
$ws.Range("G39").Value = "19"
$ws.Range("D40").Value = "0.04396"
$ws.Range("E40").Value = "-7.92%"
$ws.Range("G40").Value = "19"
$ws.Range("D41").Value = "0.007440"
$ws.Range("E41").Value = "1.29%"
$ws.Range("G41").Value = "19"
$ws.Range("D42").Value = "0.1319"
$ws.Range("E42").Value = "-3.62%"
$ws.Range("G42").Value = "19"
$ws.Range("D43").Value = "0.002114"
$ws.Range("E43").Value = "-9.54%"
$ws.Range("G43").Value = "19"
$ws.Range("D44").Value = "0.01118"
$ws.Range("E44").Value = "12.49%"
$ws.Range("G44").Value = "19"
$ws.Range("D45").Value = "0.00005974"
$ws.Range("E45").Value = "-1.28%"
$ws.Range("G45").Value = "19"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.32%"
$ws.Range("G46").Value = "19"
$ws.Range("D47").Value = "2.143"
$ws.Range("E47").Value = "161.21%"
$ws.Range("G47").Value = "19"
$ws.Range("D48").Value = "0.002403"
$ws.Range("E48").Value = "-29.20%"
$ws.Range("G48").Value = "19"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").Value = "0.32%"
$ws.Range("G49").Value = "19"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").Value = "0.32%"
$ws.Range("G50").Value = "19"
$ws.Range("G51").Value = "19"
